# Update cryptos list: price (D) and volume-1h (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.146.27"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.083.12"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +8.77%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.80"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.658"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.57%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.09"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +6.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.89"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.374"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0746"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.09%  "

$ws.Range("E12").Value = "  +6.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.389.38"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.835"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.084.23"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +8.86%  "

$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.024.86"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.53"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0827"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.56"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.23"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.42%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.35"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.03"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +12.43%  "

$ws.Range("E29").Value = "  -6.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.43"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +48.32%  "

$ws.Range("E31").Value = "  -4.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.09"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +23.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.50"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0609"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0933"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.73%  "

$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +16.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.11"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.32"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -9.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.70"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0224"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.15"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.88"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.49%  "

$ws.Range("E45").Value = "  -1.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0868"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.318.03"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.35%  "

$ws.Range("E48").Value = "  +4.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +8.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.273.14"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.23%  "

$ws.Range("E51").Value = "  -5.50%  "
